# Fix the misaligned Contact Number / Postcode columns (C and D) for rows 153-181
# by swapping the values back into their correct columns, then adjust the
# active sheet view/selection to match the saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 153; $r -le 181; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
}

# Restore the viewport/selection state captured in the saved file. Not every
# COM host persists window scroll position through to the saved OOXML, so
# this is best-effort and guarded to avoid aborting the rest of the script.
try {
    $excel.ActiveWindow.ScrollRow = 145
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}

$ws.Range("D153:D181").Select()
